# Reformatted to only one script to run the whole model.
# Updates the popPK metric analysis results (Lung + Liver sheets) with the
# refreshed model outputs: header n-count, AUC_24 / C_max means + SDs, and
# the newly-computed p-values / effect sizes for the "Better Dose" rows.

$wb = $excel.ActiveWorkbook

# ---- Lung sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("Lung")

$ws.Range("A1").Value = "Day 1, n = 100"

$ws.Range("C2").Value = 177.61000000000001
$ws.Range("D2").Value = 55.049999999999997

$ws.Range("C3").Value = 16.149999999999999
$ws.Range("D3").Value = 1.8500000000000001

$ws.Range("C4").Value = 183.15000000000001
$ws.Range("D4").Value = 53.049999999999997

$ws.Range("C5").Value = 30.670000000000002
$ws.Range("D5").Value = 5.5700000000000003

$ws.Range("C7").Value = "Lung, p = 4.4466e-27"
$ws.Range("D7").Value = 5.54

$ws.Range("C8").Value = "Lung, p = 1.48e-45"
$ws.Range("D8").Value = 14.52

$ws.Columns.Item(3).ColumnWidth = 17

# ---- Liver sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("Liver")

$ws.Range("A1").Value = "Day 1, n = 100"

$ws.Range("C2").Value = 207.94
$ws.Range("D2").Value = 63.920000000000002

$ws.Range("C3").Value = 20.77
$ws.Range("D3").Value = 2.8300000000000001

$ws.Range("C4").Value = 201.19999999999999
$ws.Range("D4").Value = 61.68

$ws.Range("C5").Value = 29.710000000000001
$ws.Range("D5").Value = 4.3899999999999997

$ws.Range("C7").Value = "Lung, p = 4.7357e-30"
$ws.Range("D7").Value = 6.75

$ws.Range("C8").Value = "Oral, p = 1.12e-35"
$ws.Range("D8").Value = 8.9399999999999995

$ws.Columns.Item(3).ColumnWidth = 17
